$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in A6 from "awesome" to "awesomeness"
$ws.Range("A6").Value = "awesomeness"

# Update the active selection to A6 (as in the diff's sheetView/selection)
$ws.Range("A6").Select()
